# Consolidate multi-run paragraphs ("A" + " " + "slide", etc.) into a
# single text run per paragraph, mirroring the PowerPoint writer change
# described in the commit ("consolidate text runs when possible").
#
# Simply re-assigning TextRange.Text to its current value is a no-op in
# this engine when the concatenated text is unchanged, so we first set a
# throwaway value and then the real target text; that forces the engine
# to rebuild the paragraph as a single run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

foreach ($shape in $s.Shapes) {
    if (-not $shape.HasTextFrame) { continue }
    $tf = $shape.TextFrame
    if (-not $tf.HasText) { continue }

    $target = $tf.TextRange.Text

    # Force a real text-replace so the writer re-emits a single <a:r>.
    $tf.TextRange.Text = "."
    $tf.TextRange.Text = $target
}
